$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "43.094.59"
$ws.Range("E2").Value = "  -2.15%  "
Set-TextValue "D3" "2.242.20"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "230.70"
$ws.Range("E5").Value = "  +0.35%  "
Set-TextValue "D6" "0.639"
$ws.Range("E6").Value = "  +1.28%  "
Set-TextValue "D7" "64.41"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue "D9" "0.437"
$ws.Range("E9").Value = "  +0.25%  "
Set-TextValue "D10" "0.0946"
$ws.Range("E10").Value = "  -6.05%  "
Set-TextValue "D11" "56.37"
$ws.Range("E11").Value = "  +0.00%  "
Set-TextValue "D12" "26.55"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("E13").Value = "  -1.59%  "
Set-TextValue "D14" "2.574.71"
$ws.Range("E14").Value = "  -0.75%  "
Set-TextValue "D15" "14.99"
$ws.Range("E15").Value = "  -4.05%  "
Set-TextValue "D16" "6.01"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("E17").Value = "  -0.13%  "
Set-TextValue "D18" "2.237.31"
$ws.Range("E18").Value = "  -0.90%  "
Set-TextValue "D19" "42.997.20"
$ws.Range("E19").Value = "  -2.03%  "
Set-TextValue "D20" "0.0₃0954"
$ws.Range("E20").Value = "  -5.44%  "
Set-TextValue "D21" "72.98"
$ws.Range("E21").Value = "  -0.68%  "
Set-TextValue "D22" "6.05"
$ws.Range("E22").Value = "  +0.78%  "
Set-TextValue "D23" "245.67"
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("E24").Value = "  +0.03%  "
Set-TextValue "D25" "3.74"
$ws.Range("E25").Value = "  +19.58%  "
Set-TextValue "D26" "2.42"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -2.05%  "
Set-TextValue "D28" "173.98"
$ws.Range("E28").Value = "  +1.07%  "
Set-TextValue "D29" "9.67"
$ws.Range("E29").Value = "  -2.89%  "
Set-TextValue "D30" "21.53"
$ws.Range("E30").Value = "  +3.80%  "
Set-TextValue "D31" "1.41"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("E33").Value = "  +0.70%  "
Set-TextValue "D34" "4.91"
$ws.Range("E34").Value = "  +4.65%  "
Set-TextValue "D35" "0.0676"
$ws.Range("E35").Value = "  -0.27%  "
Set-TextValue "D36" "4.88"
$ws.Range("E36").Value = "  +0.01%  "
Set-TextValue "D37" "3.59"
$ws.Range("E37").Value = "  -6.74%  "
Set-TextValue "D38" "6.28"
$ws.Range("E38").Value = "  -6.36%  "
Set-TextValue "D39" "2.26"
$ws.Range("E39").Value = "  -2.44%  "
Set-TextValue "D40" "0.0247"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -0.05%  "
Set-TextValue "D42" "8.67"
$ws.Range("E42").Value = "  +3.86%  "
Set-TextValue "D43" "4.43"
$ws.Range("E43").Value = "  +0.36%  "
Set-TextValue "D44" "17.02"
$ws.Range("E44").Value = "  -2.28%  "
Set-TextValue "D45" "96.39"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D46" "1.18"
$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D47" "0.0934"
$ws.Range("E47").Value = "  -2.31%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D48" "1.433.86"
$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("B49").Value = "TerraClassic"
$ws.Range("C49").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue "D49" "0.000205"
$ws.Range("E49").Value = "  -0.66%  "

Set-TextValue "D50" "9.87"
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("E51").Value = "  -1.41%  "
